$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "now价格" column (C) stores its values as text (e.g. "6.30", "10.80")
# so that formatted strings like trailing zeros are preserved exactly.
# Mark the range as Text before writing so the new values keep that type.
$ws.Range("C2:C9").NumberFormat = "@"

$ws.Range("C2").Value = "6.43"
$ws.Range("C3").Value = "72.09"
$ws.Range("C4").Value = "13.31"
$ws.Range("C5").Value = "17.40"
$ws.Range("C6").Value = "42.39"
$ws.Range("C7").Value = "1488.11"
$ws.Range("C8").Value = "360.67"
$ws.Range("C9").Value = "10.75"
